$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new worksheet named "2022-Q3" right before the existing
#    "2022-Q2" sheet (i.e. as the 2nd sheet, right after "总计").
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($refSheet)
$q3.Name = "2022-Q3"

# Header row (B1:H1), bold + centered + boxed, matching the look of the
# header rows used on every other quarter sheet.
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$headerRng = $q3.Range("B1:H1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

# Data rows. Column A is a plain numeric row index (0-based). Columns D, E,
# F and G are textual percentages/amounts (kept as text so formatting such
# as trailing zeros is preserved exactly); column H is a plain number.
$rows = @(
    @(0, "006551", "中庚价值领航混合",              "118.19", "91.86", "3.46", "4.0894", 8),
    @(1, "007130", "中庚小盘价值股票",              "75.87",  "93.06", "4.92", "3.7328", 4),
    @(2, "007497", "中庚价值灵动灵活配置混合",        "36.46",  "89.30", "2.11", "0.7693", 8),
    @(3, "001551", "天弘中证医药100指数型发起式 C",  "8.58",   "95.24", "1.59", "0.1364", 1),
    @(4, "001550", "天弘中证医药100指数型发起式 A",  "5.31",   "95.24", "1.59", "0.0844", 1),
    @(5, "210011", "金鹰灵活配置混合C",              "1.86",   "24.78", "0.54", "0.0100", 7),
    @(6, "210010", "金鹰灵活配置混合A",              "1.46",   "24.78", "0.54", "0.0079", 7),
    @(7, "005264", "国都多策略混合",                "0.04",   "56.56", "3.14", "0.0013", 9)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $aCell = $q3.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    $bCell = $q3.Cells.Item($r, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $row[1]

    $q3.Cells.Item($r, 3).Value = $row[2]

    $dCell = $q3.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[3]

    $eCell = $q3.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[4]

    $fCell = $q3.Cells.Item($r, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $row[5]

    $gCell = $q3.Cells.Item($r, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $row[6]

    $q3.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------------
# 2) Update the "总计" summary sheet: push the existing quarterly rows down
#    by one and insert the brand-new 2022-Q3 summary row at the top.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

for ($r = 8; $r -ge 2; $r--) {
    $dest = $r + 1
    $total.Cells.Item($r, 1).Copy($total.Cells.Item($dest, 1))
    $total.Cells.Item($dest, 1).Value = $dest - 2
    $total.Cells.Item($dest, 2).Value = $total.Cells.Item($r, 2).Value()
    $total.Cells.Item($dest, 3).Value = $total.Cells.Item($r, 3).Value()
    $total.Cells.Item($dest, 4).Value = $total.Cells.Item($r, 4).Value()
}

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 8
$total.Cells.Item(2, 4).Value = 8.83

# Restore the original active sheet/selection (adding a sheet shifts focus
# onto it) so the workbook-level view state is left as it was.
[void]$total.Activate()
[void]$total.Range("A1").Select()
